$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header updates: report number and date range (new week) ---
$ws.Range("A8").Value = "Volume 30   Number  18"
$ws.Range("C9").Value = "Report Covering the Week  5/1/2023  Through  5/7/2023"

# --- Weekly crime complaint stats table (rows 14-30): refreshed figures ---
$ws.Range("D14").NumberFormat = "#,##0"
$ws.Range("D14").Value = 1
$ws.Range("E14").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E14").Value = -100
$ws.Range("J14").Value = 4
$ws.Range("K14").Value = -75
$ws.Range("N14").Value = -95.238095238095
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 1
$ws.Range("F15").NumberFormat = "#,##0"
$ws.Range("F15").Value = 2
$ws.Range("I15").Value = 8
$ws.Range("K15").Value = 33.333333333333
$ws.Range("L15").Value = 14.285714285714
$ws.Range("M15").Value = 166.666666666667
$ws.Range("N15").Value = -71.428571428571
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 11
$ws.Range("E16").Value = -72.727272727272
$ws.Range("F16").Value = 22
$ws.Range("G16").Value = 35
$ws.Range("H16").Value = -37.142857142857
$ws.Range("I16").Value = 105
$ws.Range("J16").Value = 134
$ws.Range("K16").Value = -21.641791044776
$ws.Range("L16").Value = -11.764705882352
$ws.Range("M16").Value = -27.083333333333
$ws.Range("N16").Value = -86.725663716814
$ws.Range("C17").Value = 13
$ws.Range("D17").Value = 17
$ws.Range("E17").Value = -23.529411764705
$ws.Range("G17").Value = 52
$ws.Range("H17").Value = -9.615384615384
$ws.Range("I17").Value = 222
$ws.Range("J17").Value = 231
$ws.Range("K17").Value = -3.896103896103
$ws.Range("L17").Value = 16.842105263157
$ws.Range("M17").Value = 64.444444444444
$ws.Range("N17").Value = -37.110481586402
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -35.714285714285
$ws.Range("I18").Value = 57
$ws.Range("J18").Value = 72
$ws.Range("K18").Value = -20.833333333333
$ws.Range("L18").Value = -8.064516129032
$ws.Range("M18").Value = -18.571428571428
$ws.Range("N18").Value = -92.359249329758
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -33.333333333333
$ws.Range("G19").Value = 41
$ws.Range("H19").Value = -26.829268292682
$ws.Range("I19").Value = 173
$ws.Range("J19").Value = 216
$ws.Range("K19").Value = -19.907407407407
$ws.Range("L19").Value = -10.824742268041
$ws.Range("M19").Value = 90.10989010989
$ws.Range("N19").Value = -43.831168831168
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 19
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 111.111111111111
$ws.Range("I20").Value = 100
$ws.Range("J20").Value = 83
$ws.Range("K20").Value = 20.481927710843
$ws.Range("L20").Value = 170.27027027027
$ws.Range("M20").Value = 156.410256410256
$ws.Range("N20").Value = -71.988795518207
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 50
$ws.Range("E21").Value = -38
$ws.Range("F21").Value = 129
$ws.Range("G21").Value = 152
$ws.Range("H21").Value = -15.131578947368
$ws.Range("I21").Value = 666
$ws.Range("J21").Value = 746
$ws.Range("K21").Value = -10.723860589812
$ws.Range("L21").Value = 9.001636661211
$ws.Range("M21").Value = 37.603305785124
$ws.Range("N21").Value = -74.42396313364
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("C22").Value = 1
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Value = "0"
$ws.Range("E22").NumberFormat = "General"
$ws.Range("E22").Value = "***.*"
$ws.Range("F22").NumberFormat = "#,##0"
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 6
$ws.Range("K22").Value = -45.454545454545
$ws.Range("L22").Value = -60
$ws.Range("M22").Value = -25
$ws.Range("D23").NumberFormat = "#,##0"
$ws.Range("D23").Value = 1
$ws.Range("E23").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 2
$ws.Range("G23").NumberFormat = "#,##0"
$ws.Range("G23").Value = 1
$ws.Range("H23").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("H23").Value = 100
$ws.Range("J23").Value = 9
$ws.Range("K23").Value = 33.333333333333
$ws.Range("L23").Value = 20
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -26.666666666666
$ws.Range("F24").Value = 83
$ws.Range("G24").Value = 105
$ws.Range("H24").Value = -20.952380952381
$ws.Range("I24").Value = 326
$ws.Range("J24").Value = 458
$ws.Range("K24").Value = -28.82096069869
$ws.Range("L24").Value = 10.884353741496
$ws.Range("M24").Value = -2.395209580838
$ws.Range("C25").Value = 27
$ws.Range("D25").Value = 22
$ws.Range("E25").Value = 22.727272727272
$ws.Range("F25").Value = 111
$ws.Range("G25").Value = 78
$ws.Range("H25").Value = 42.307692307692
$ws.Range("I25").Value = 426
$ws.Range("J25").Value = 340
$ws.Range("K25").Value = 25.294117647058
$ws.Range("L25").Value = 56.043956043956
$ws.Range("M25").Value = 1.670644391408
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("C26").Value = 2
$ws.Range("E26").Value = 100
$ws.Range("F26").NumberFormat = "#,##0"
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 19
$ws.Range("J26").Value = 15
$ws.Range("K26").Value = 26.666666666666
$ws.Range("L26").Value = 35.714285714285
$ws.Range("C27").NumberFormat = "General"
$ws.Range("C27").Value = "0"
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -57.142857142857
$ws.Range("I27").Value = 22
$ws.Range("J27").Value = 31
$ws.Range("K27").Value = -29.032258064516
$ws.Range("L27").Value = -24.137931034482
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("C28").Value = 2
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("D28").Value = 3
$ws.Range("E28").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E28").Value = -33.333333333333
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 33.333333333333
$ws.Range("I28").Value = 10
$ws.Range("J28").Value = 13
$ws.Range("K28").Value = -23.076923076923
$ws.Range("L28").Value = 25
$ws.Range("M28").Value = -33.333333333333
$ws.Range("N28").Value = -81.481481481481
$ws.Range("C29").NumberFormat = "#,##0"
$ws.Range("C29").Value = 1
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("D29").Value = 1
$ws.Range("E29").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 7
$ws.Range("J29").Value = 11
$ws.Range("K29").Value = -36.363636363636
$ws.Range("L29").Value = -12.5
$ws.Range("M29").Value = -50
$ws.Range("N29").Value = -86.274509803921
